# "Added browser as a command option"
#
# The QBE test fixture's "Status" column (AI) previously recorded
# "Passed" for test cases 2-6. With a browser command option now
# added to the test runner, those recorded results are stale and are
# cleared back out to match the still-unset rows below them (AI7:AI9).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("AI2:AI6").Value = ""
